$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Make Attendance" action now just stamps the date header; the
# previously-recorded attendance rows are cleared so the (new) download
# button can populate them fresh.
# Format A1 as Text first so the date-like string "2021/1/19" is stored
# as literal text rather than being auto-converted to a date serial.
$ws.Range("A1").NumberFormat = "@"
$ws.Range("A1").Value = "2021/1/19"

# Remove the old attendance entries (rows 3-4), leaving just the two
# header rows.
$ws.Rows("3:4").Delete()
